$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.792.53"
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").Value = "1.916.23"
$ws.Range("E3").Value = "  +1.55%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.96%  "

$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4915"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2981"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06778"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").Value = "1.923.39"
$ws.Range("E10").Value = "  +1.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07375"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.207"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.96%  "

$ws.Range("E14").Value = "  -1.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6746"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.24%  "

$ws.Range("D16").Value = "30.771.93"
$ws.Range("E16").Value = "  +0.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008004"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.02%  "

$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("D20").Value = "2.147.47"
$ws.Range("E20").Value = "  +0.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.278"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "203.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.300"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.670"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.977"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.439"
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.366"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09217"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.53%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.087"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.70%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05441"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7519"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.80%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.123"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.701"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01869"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.729"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9295"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.083"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.80%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4509"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.37%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +26.57%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "107.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.940"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.003"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1394"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.743"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "36.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.134"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05986"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4068"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.99%  "
